# Update "想去人数" (want-to-go count) figures on the 苏州-漫展信息 workbook.
# Sheet "展览" and sheet "全部类型" share most rows; sheet "演出" has one row
# in common. Values below were re-scraped and bumped slightly upward.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 551
$ws1.Range("F4").Value = 1571
$ws1.Range("F8").Value = 185
$ws1.Range("F9").Value = 761
$ws1.Range("F13").Value = 65
$ws1.Range("F14").Value = 516
$ws1.Range("F15").Value = 24
$ws1.Range("F16").Value = 6539
$ws1.Range("F22").Value = 7
$ws1.Range("F23").Value = 15630
$ws1.Range("F24").Value = 1544
$ws1.Range("F25").Value = 14
$ws1.Range("F26").Value = 302
$ws1.Range("F27").Value = 154
$ws1.Range("F28").Value = 107
$ws1.Range("F29").Value = 11131
$ws1.Range("F30").Value = 788
$ws1.Range("F32").Value = 263
$ws1.Range("F33").Value = 379

# --- Sheet: 演出 (performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 43

# --- Sheet: 全部类型 (all types, combined listing) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 551
$ws4.Range("F4").Value = 1571
$ws4.Range("F9").Value = 185
$ws4.Range("F10").Value = 761
$ws4.Range("F15").Value = 65
$ws4.Range("F16").Value = 516
$ws4.Range("F17").Value = 43
$ws4.Range("F18").Value = 24
$ws4.Range("F19").Value = 6539
$ws4.Range("F26").Value = 7
$ws4.Range("F27").Value = 15631
$ws4.Range("F28").Value = 1544
$ws4.Range("F29").Value = 14
$ws4.Range("F30").Value = 302
$ws4.Range("F31").Value = 154
$ws4.Range("F32").Value = 107
$ws4.Range("F34").Value = 11131
$ws4.Range("F35").Value = 788
$ws4.Range("F37").Value = 263
$ws4.Range("F38").Value = 379
